$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    # Force the cell to Text format so Excel does not coerce numeric-
    # looking strings (e.g. "0.603", "580.70") into numbers, then
    # restore the default "Normal" style so no stray style index is
    # left behind on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '67.948.12'
$ws.Range('E2').Value = '  -0.97%  '
# Row 3
Set-TextValue $ws 'D3' '3.269.41'
$ws.Range('E3').Value = '  +0.14%  '
# Row 4
$ws.Range('E4').Value = '  -0.06%  '
# Row 5
Set-TextValue $ws 'D5' '580.70'
$ws.Range('E5').Value = '  -0.46%  '
# Row 6
Set-TextValue $ws 'D6' '183.60'
$ws.Range('E6').Value = '  +1.21%  '
# Row 7
$ws.Range('E7').Value = '  +0.01%  '
# Row 8
Set-TextValue $ws 'D8' '0.603'
$ws.Range('E8').Value = '  +1.15%  '
# Row 9
$ws.Range('E9').Value = '  -2.09%  '
# Row 10
$ws.Range('E10').Value = '  -1.29%  '
# Row 11
Set-TextValue $ws 'D11' '0.409'
$ws.Range('E11').Value = '  -3.49%  '
# Row 12
Set-TextValue $ws 'D12' '3.836.12'
$ws.Range('E12').Value = '  +0.04%  '
# Row 13
$ws.Range('E13').Value = '  +0.57%  '
# Row 14
Set-TextValue $ws 'D14' '27.42'
$ws.Range('E14').Value = '  -3.68%  '
# Row 15
Set-TextValue $ws 'D15' '67.915.81'
$ws.Range('E15').Value = '  -1.00%  '
# Row 16
$ws.Range('E16').Value = '  -1.70%  '
# Row 17
Set-TextValue $ws 'D17' '3.270.97'
$ws.Range('E17').Value = '  +1.67%  '
# Row 18
$ws.Range('E18').Value = '  -2.09%  '
# Row 19
Set-TextValue $ws 'D19' '13.39'
$ws.Range('E19').Value = '  -1.29%  '
# Row 20
Set-TextValue $ws 'D20' '403.94'
$ws.Range('E20').Value = '  +2.53%  '
# Row 21
$ws.Range('E21').Value = '  -1.97%  '
# Row 22
$ws.Range('E22').Value = '  +0.00%  '
# Row 23
Set-TextValue $ws 'D23' '71.12'
$ws.Range('E23').Value = '  -1.29%  '
# Row 24
$ws.Range('E24').Value = '  -1.46%  '
# Row 25
$ws.Range('E25').Value = '  -1.58%  '
# Row 26
$ws.Range('E26').Value = '  -0.67%  '
# Row 27
$ws.Range('E27').Value = '  -1.42%  '
# Row 28
$ws.Range('E28').Value = '  +0.38%  '
# Row 29
Set-TextValue $ws 'D29' '1.95'
$ws.Range('E29').Value = '  -1.49%  '
# Row 30
Set-TextValue $ws 'D30' '22.70'
$ws.Range('E30').Value = '  -1.06%  '
# Row 31
Set-TextValue $ws 'D31' '5.47'
$ws.Range('E31').Value = '  -3.81%  '
# Row 32
Set-TextValue $ws 'D32' '6.90'
$ws.Range('E32').Value = '  -3.24%  '
# Row 33
$ws.Range('E33').Value = '  +0.07%  '
# Row 34
Set-TextValue $ws 'D34' '1.25'
$ws.Range('E34').Value = '  -2.77%  '
# Row 35
Set-TextValue $ws 'D35' '164.36'
$ws.Range('E35').Value = '  +0.19%  '
# Row 36
$ws.Range('E36').Value = '  -3.14%  '
# Row 37
$ws.Range('E37').Value = '  -1.30%  '
# Row 38
Set-TextValue $ws 'D38' '27.12'
$ws.Range('E38').Value = '  +2.81%  '
# Row 39
Set-TextValue $ws 'D39' '0.801'
$ws.Range('E39').Value = '  -3.43%  '
# Row 40
Set-TextValue $ws 'D40' '4.49'
$ws.Range('E40').Value = '  -2.22%  '
# Row 41
Set-TextValue $ws 'D41' '6.35'
$ws.Range('E41').Value = '  -3.26%  '
# Row 42
Set-TextValue $ws 'D42' '2.678.59'
$ws.Range('E42').Value = '  +2.62%  '
# Row 43
$ws.Range('E43').Value = '  -1.22%  '
# Row 44
Set-TextValue $ws 'D44' '0.0678'
$ws.Range('E44').Value = '  -1.43%  '
# Row 45
$ws.Range('E45').Value = '  -2.37%  '
# Row 46
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D46' '24.61'
$ws.Range('E46').Value = '  +0.02%  '
# Row 47
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D47' '335.11'
$ws.Range('E47').Value = '  -3.67%  '
# Row 48
Set-TextValue $ws 'D48' '0.0275'
$ws.Range('E48').Value = '  -2.47%  '
# Row 49
Set-TextValue $ws 'D49' '6.30'
# Row 50
$ws.Range('E50').Value = '  -1.26%  '
# Row 51
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws 'D51' '30.78'
$ws.Range('E51').Value = '  -2.60%  '
